$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 112
$ws.Range("F6").Value = 525
$ws.Range("F7").Value = 4923
$ws.Range("F8").Value = 4923
$ws.Range("F9").Value = 115
$ws.Range("F10").Value = 132
$ws.Range("F11").Value = 185
$ws.Range("F12").Value = 33
$ws.Range("F14").Value = 143
$ws.Range("F15").Value = 7973
$ws.Range("F16").Value = 7973
$ws.Range("F17").Value = 262
$ws.Range("F19").Value = 567
$ws.Range("F20").Value = 2464
$ws.Range("F21").Value = 6312
$ws.Range("F22").Value = 2278
$ws.Range("F25").Value = 2504
$ws.Range("F27").Value = 5
$ws.Range("F28").Value = 6294
$ws.Range("F29").Value = 167
$ws.Range("F30").Value = 53
$ws.Range("F31").Value = 123
$ws.Range("F32").Value = 97
$ws.Range("F34").Value = 6675
$ws.Range("F39").Value = 8
$ws.Range("F42").Value = 32
$ws.Range("F43").Value = 2491
$ws.Range("F47").Value = 52
$ws.Range("F48").Value = 473
$ws.Range("F49").Value = 2185
$ws.Range("F50").Value = 61
$ws.Range("F51").Value = 1101
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 142
$ws.Range("F10").Value = 46
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1466
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1466
$ws.Range("F7").Value = 112
$ws.Range("F9").Value = 525
$ws.Range("F10").Value = 4923
$ws.Range("F11").Value = 4923
$ws.Range("F12").Value = 115
$ws.Range("F13").Value = 132
$ws.Range("F14").Value = 185
$ws.Range("F15").Value = 33
$ws.Range("F16").Value = 143
$ws.Range("F17").Value = 7973
$ws.Range("F18").Value = 7973
$ws.Range("F19").Value = 262
$ws.Range("F21").Value = 567
$ws.Range("F22").Value = 2464
$ws.Range("F24").Value = 142
$ws.Range("F25").Value = 6312
$ws.Range("F26").Value = 2278
$ws.Range("F27").Value = 2504
$ws.Range("F29").Value = 46
$ws.Range("F30").Value = 5
$ws.Range("F31").Value = 6294
$ws.Range("F32").Value = 167
$ws.Range("F33").Value = 53
$ws.Range("F34").Value = 123
$ws.Range("F35").Value = 97
$ws.Range("F37").Value = 6675
$ws.Range("F41").Value = 32
$ws.Range("F43").Value = 2491
$ws.Range("F46").Value = 52
$ws.Range("F47").Value = 473
$ws.Range("F49").Value = 2185
$ws.Range("F50").Value = 61
$ws.Range("F51").Value = 1101
